$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ALC_updates = @(
    @{ Row = 8; Col = "H"; Value = 2750 },
    @{ Row = 8; Col = "I"; Value = 250.5 },
    @{ Row = 8; Col = "J"; Value = 5249.5 },
    @{ Row = 8; Col = "K"; Value = 751.5 },
    @{ Row = 8; Col = "L"; Value = 15748.5 },
    @{ Row = 8; Col = "M"; Value = -612.5 },
    @{ Row = 8; Col = "N"; Value = -16026.5 },
    @{ Row = 9; Col = "H"; Value = 87.5 },
    @{ Row = 9; Col = "I"; Value = 87.5 },
    @{ Row = 9; Col = "J"; Value = 0 },
    @{ Row = 9; Col = "K"; Value = 87.5 },
    @{ Row = 9; Col = "L"; Value = 0 },
    @{ Row = 9; Col = "M"; Value = 81.5 },
    @{ Row = 11; Col = "H"; Value = 32.666668 },
    @{ Row = 11; Col = "I"; Value = 32.666668 },
    @{ Row = 11; Col = "K"; Value = 32.666668 },
    @{ Row = 11; Col = "M"; Value = 107.333332 },
    @{ Row = 74; Col = "H"; Value = 3 },
    @{ Row = 74; Col = "I"; Value = 3 },
    @{ Row = 74; Col = "K"; Value = 3 },
    @{ Row = 74; Col = "M"; Value = 933 },
    @{ Row = 77; Col = "H"; Value = 3 },
    @{ Row = 77; Col = "I"; Value = 3 },
    @{ Row = 77; Col = "K"; Value = 15 },
    @{ Row = 77; Col = "M"; Value = 4665 },
    @{ Row = 98; Col = "H"; Value = 13499.889 },
    @{ Row = 98; Col = "J"; Value = 13499.889 },
    @{ Row = 98; Col = "L"; Value = 13499.889 },
    @{ Row = 98; Col = "N"; Value = -16495.889 },
    @{ Row = 100; Col = "H"; Value = 1675 },
    @{ Row = 100; Col = "I"; Value = 2200 },
    @{ Row = 100; Col = "J"; Value = 1500 },
    @{ Row = 100; Col = "K"; Value = 2200 },
    @{ Row = 100; Col = "L"; Value = 1500 },
    @{ Row = 100; Col = "M"; Value = -1659 },
    @{ Row = 100; Col = "N"; Value = -2582 },
    @{ Row = 122; Col = "H"; Value = 13499.889 },
    @{ Row = 122; Col = "J"; Value = 13499.889 },
    @{ Row = 122; Col = "L"; Value = 40499.667 },
    @{ Row = 122; Col = "N"; Value = -45399.667 },
    @{ Row = 132; Col = "H"; Value = 6799.4 },
    @{ Row = 132; Col = "I"; Value = 4998.8 },
    @{ Row = 132; Col = "K"; Value = 14996.4 },
    @{ Row = 132; Col = "M"; Value = -12466.4 }
)
foreach ($u in $ALC_updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Value
}
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ARM_updates = @(
    @{ Row = 32; Col = "H"; Value = 1618 },
    @{ Row = 32; Col = "I"; Value = 1502.8334 },
    @{ Row = 32; Col = "K"; Value = 1502.8334 },
    @{ Row = 32; Col = "M"; Value = -1215.8334 },
    @{ Row = 63; Col = "H"; Value = 2557 },
    @{ Row = 63; Col = "I"; Value = 2724.25 },
    @{ Row = 63; Col = "J"; Value = 1888 },
    @{ Row = 63; Col = "K"; Value = 2724.25 },
    @{ Row = 63; Col = "L"; Value = 1888 },
    @{ Row = 63; Col = "M"; Value = -2038.25 },
    @{ Row = 63; Col = "N"; Value = -3260 },
    @{ Row = 66; Col = "H"; Value = 2557 },
    @{ Row = 66; Col = "I"; Value = 2724.25 },
    @{ Row = 66; Col = "J"; Value = 1888 },
    @{ Row = 66; Col = "K"; Value = 13621.25 },
    @{ Row = 66; Col = "L"; Value = 9440 },
    @{ Row = 66; Col = "M"; Value = -10189.25 },
    @{ Row = 66; Col = "N"; Value = -16304 },
    @{ Row = 104; Col = "H"; Value = 17999.5 },
    @{ Row = 104; Col = "J"; Value = 17999.5 },
    @{ Row = 104; Col = "L"; Value = 17999.5 },
    @{ Row = 104; Col = "N"; Value = -24987.5 },
    @{ Row = 112; Col = "H"; Value = 29999.5 },
    @{ Row = 112; Col = "J"; Value = 29999.5 },
    @{ Row = 112; Col = "L"; Value = 29999.5 },
    @{ Row = 112; Col = "N"; Value = -32953.5 },
    @{ Row = 114; Col = "H"; Value = 49999 },
    @{ Row = 114; Col = "J"; Value = 49999 },
    @{ Row = 114; Col = "L"; Value = 49999 },
    @{ Row = 114; Col = "N"; Value = -58677 }
)
foreach ($u in $ARM_updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Value
}

$ws = $wb.Worksheets.Item("BSM")
$BSM_updates = @(
    @{ Row = 11; Col = "H"; Value = 449.25 },
    @{ Row = 11; Col = "I"; Value = 266 },
    @{ Row = 11; Col = "J"; Value = 999 },
    @{ Row = 11; Col = "K"; Value = 266 },
    @{ Row = 11; Col = "L"; Value = 999 },
    @{ Row = 11; Col = "M"; Value = -126 },
    @{ Row = 11; Col = "N"; Value = -1279 },
    @{ Row = 36; Col = "H"; Value = 6894 },
    @{ Row = 36; Col = "I"; Value = 6894 },
    @{ Row = 36; Col = "K"; Value = 6894 },
    @{ Row = 36; Col = "M"; Value = -6360 },
    @{ Row = 87; Col = "H"; Value = 50321 },
    @{ Row = 87; Col = "I"; Value = 50321 },
    @{ Row = 87; Col = "J"; Value = 0 },
    @{ Row = 87; Col = "K"; Value = 50321 },
    @{ Row = 87; Col = "L"; Value = 0 },
    @{ Row = 87; Col = "M"; Value = -49073 },
    @{ Row = 90; Col = "H"; Value = 50321 },
    @{ Row = 90; Col = "I"; Value = 50321 },
    @{ Row = 90; Col = "J"; Value = 0 },
    @{ Row = 90; Col = "K"; Value = 150963 },
    @{ Row = 90; Col = "L"; Value = 0 },
    @{ Row = 90; Col = "M"; Value = -144723 },
    @{ Row = 94; Col = "H"; Value = 2654.5 },
    @{ Row = 94; Col = "I"; Value = 2654.5 },
    @{ Row = 94; Col = "K"; Value = 2654.5 },
    @{ Row = 94; Col = "M"; Value = -2203.5 }
)
foreach ($u in $BSM_updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Value
}
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$CRP_updates = @(
    @{ Row = 4; Col = "H"; Value = 12502100 },
    @{ Row = 4; Col = "J"; Value = 12502100 },
    @{ Row = 4; Col = "L"; Value = 12502100 },
    @{ Row = 4; Col = "N"; Value = -12502324 },
    @{ Row = 8; Col = "H"; Value = 249 },
    @{ Row = 8; Col = "I"; Value = 249 },
    @{ Row = 8; Col = "K"; Value = 249 },
    @{ Row = 8; Col = "M"; Value = -109 }
)
foreach ($u in $CRP_updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Value
}

$ws = $wb.Worksheets.Item("CUL")
$CUL_updates = @(
    @{ Row = 2; Col = "H"; Value = 30.875 },
    @{ Row = 2; Col = "I"; Value = 36.076923 },
    @{ Row = 2; Col = "J"; Value = 8.333333 },
    @{ Row = 2; Col = "K"; Value = 216.461538 },
    @{ Row = 2; Col = "L"; Value = 49.999998 },
    @{ Row = 2; Col = "M"; Value = -103.461538 },
    @{ Row = 2; Col = "N"; Value = -275.999998 },
    @{ Row = 4; Col = "H"; Value = 1115457.2 },
    @{ Row = 4; Col = "I"; Value = 269.3846 },
    @{ Row = 4; Col = "J"; Value = 2150989 },
    @{ Row = 4; Col = "K"; Value = 808.1537999999999 },
    @{ Row = 4; Col = "L"; Value = 6452967 },
    @{ Row = 4; Col = "M"; Value = -696.1537999999999 },
    @{ Row = 4; Col = "N"; Value = -6453191 },
    @{ Row = 5; Col = "H"; Value = 1382.7333 },
    @{ Row = 5; Col = "J"; Value = 1426.5834 },
    @{ Row = 5; Col = "L"; Value = 4279.7502 },
    @{ Row = 5; Col = "N"; Value = -4503.7502 },
    @{ Row = 45; Col = "H"; Value = 0 },
    @{ Row = 45; Col = "J"; Value = 0 },
    @{ Row = 45; Col = "L"; Value = 0 },
    @{ Row = 81; Col = "H"; Value = 0 },
    @{ Row = 81; Col = "I"; Value = 0 },
    @{ Row = 81; Col = "J"; Value = 0 },
    @{ Row = 81; Col = "K"; Value = 0 },
    @{ Row = 81; Col = "L"; Value = 0 },
    @{ Row = 84; Col = "H"; Value = 0 },
    @{ Row = 84; Col = "I"; Value = 0 },
    @{ Row = 84; Col = "J"; Value = 0 },
    @{ Row = 84; Col = "K"; Value = 0 },
    @{ Row = 84; Col = "L"; Value = 0 },
    @{ Row = 109; Col = "H"; Value = 0 },
    @{ Row = 109; Col = "I"; Value = 0 },
    @{ Row = 109; Col = "K"; Value = 0 },
    @{ Row = 122; Col = "H"; Value = 414 },
    @{ Row = 122; Col = "J"; Value = 849 },
    @{ Row = 122; Col = "L"; Value = 7641 },
    @{ Row = 122; Col = "N"; Value = -12541 },
    @{ Row = 129; Col = "H"; Value = 2248.8333 },
    @{ Row = 129; Col = "I"; Value = 1999 },
    @{ Row = 129; Col = "J"; Value = 2298.8 },
    @{ Row = 129; Col = "K"; Value = 5997 },
    @{ Row = 129; Col = "L"; Value = 6896.400000000001 },
    @{ Row = 129; Col = "M"; Value = -997 },
    @{ Row = 129; Col = "N"; Value = -16896.4 },
    @{ Row = 131; Col = "H"; Value = 1976.4706 },
    @{ Row = 131; Col = "I"; Value = 810.8 },
    @{ Row = 131; Col = "J"; Value = 2462.1667 },
    @{ Row = 131; Col = "K"; Value = 2432.4 },
    @{ Row = 131; Col = "L"; Value = 7386.500100000001 },
    @{ Row = 131; Col = "M"; Value = 2607.6 },
    @{ Row = 131; Col = "N"; Value = -17466.5001 },
    @{ Row = 135; Col = "H"; Value = 1382.7333 },
    @{ Row = 135; Col = "J"; Value = 1426.5834 },
    @{ Row = 135; Col = "L"; Value = 12839.2506 },
    @{ Row = 135; Col = "N"; Value = -17909.2506 }
)
foreach ($u in $CUL_updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Value
}
$ws.Range("N45").ClearContents()
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
$ws.Range("M109").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$GSM_updates = @(
    @{ Row = 52; Col = "H"; Value = 40000 },
    @{ Row = 52; Col = "J"; Value = 40000 },
    @{ Row = 52; Col = "L"; Value = 40000 },
    @{ Row = 52; Col = "N"; Value = -40518 },
    @{ Row = 80; Col = "H"; Value = 1990.75 },
    @{ Row = 80; Col = "I"; Value = 2304.3333 },
    @{ Row = 80; Col = "J"; Value = 1050 },
    @{ Row = 80; Col = "K"; Value = 2304.3333 },
    @{ Row = 80; Col = "L"; Value = 1050 },
    @{ Row = 80; Col = "M"; Value = -1306.3333 },
    @{ Row = 80; Col = "N"; Value = -3046 },
    @{ Row = 83; Col = "H"; Value = 1990.75 },
    @{ Row = 83; Col = "I"; Value = 2304.3333 },
    @{ Row = 83; Col = "J"; Value = 1050 },
    @{ Row = 83; Col = "K"; Value = 11521.6665 },
    @{ Row = 83; Col = "L"; Value = 5250 },
    @{ Row = 83; Col = "M"; Value = -6529.666499999999 },
    @{ Row = 83; Col = "N"; Value = -15234 },
    @{ Row = 92; Col = "H"; Value = 105125 },
    @{ Row = 92; Col = "J"; Value = 105125 },
    @{ Row = 92; Col = "L"; Value = 105125 },
    @{ Row = 92; Col = "N"; Value = -108869 },
    @{ Row = 94; Col = "H"; Value = 0 },
    @{ Row = 94; Col = "J"; Value = 0 },
    @{ Row = 94; Col = "L"; Value = 0 },
    @{ Row = 97; Col = "H"; Value = 300 },
    @{ Row = 97; Col = "I"; Value = 300 },
    @{ Row = 97; Col = "K"; Value = 300 },
    @{ Row = 97; Col = "M"; Value = 196 },
    @{ Row = 113; Col = "H"; Value = 3376.375 },
    @{ Row = 113; Col = "I"; Value = 3144.4285 },
    @{ Row = 113; Col = "K"; Value = 3144.4285 },
    @{ Row = 113; Col = "M"; Value = -974.4285 },
    @{ Row = 126; Col = "H"; Value = 2000 },
    @{ Row = 126; Col = "I"; Value = 2000 },
    @{ Row = 126; Col = "J"; Value = 0 },
    @{ Row = 126; Col = "K"; Value = 6000 },
    @{ Row = 126; Col = "L"; Value = 0 },
    @{ Row = 126; Col = "M"; Value = -3530 }
)
foreach ($u in $GSM_updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Value
}
$ws.Range("N94").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$LTW_updates = @(
    @{ Row = 25; Col = "H"; Value = 3000 },
    @{ Row = 25; Col = "I"; Value = 3000 },
    @{ Row = 25; Col = "K"; Value = 3000 },
    @{ Row = 25; Col = "M"; Value = -2770 },
    @{ Row = 80; Col = "H"; Value = 25000 },
    @{ Row = 80; Col = "J"; Value = 25000 },
    @{ Row = 80; Col = "L"; Value = 25000 },
    @{ Row = 80; Col = "N"; Value = -27246 },
    @{ Row = 83; Col = "H"; Value = 25000 },
    @{ Row = 83; Col = "J"; Value = 25000 },
    @{ Row = 83; Col = "L"; Value = 75000 },
    @{ Row = 83; Col = "N"; Value = -86232 },
    @{ Row = 92; Col = "H"; Value = 31353 },
    @{ Row = 92; Col = "I"; Value = 31353 },
    @{ Row = 92; Col = "K"; Value = 31353 },
    @{ Row = 92; Col = "M"; Value = -28857 },
    @{ Row = 110; Col = "H"; Value = 22322 },
    @{ Row = 110; Col = "J"; Value = 22322 },
    @{ Row = 110; Col = "L"; Value = 22322 },
    @{ Row = 110; Col = "N"; Value = -30502 },
    @{ Row = 139; Col = "H"; Value = 97977 },
    @{ Row = 139; Col = "J"; Value = 97977 },
    @{ Row = 139; Col = "L"; Value = 97977 },
    @{ Row = 139; Col = "N"; Value = -108257 }
)
foreach ($u in $LTW_updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Value
}

$ws = $wb.Worksheets.Item("WVR")
$WVR_updates = @(
    @{ Row = 105; Col = "H"; Value = 33488.668 },
    @{ Row = 105; Col = "J"; Value = 33488.668 },
    @{ Row = 105; Col = "L"; Value = 33488.668 },
    @{ Row = 105; Col = "N"; Value = -40476.668 },
    @{ Row = 136; Col = "H"; Value = 10576.556 },
    @{ Row = 136; Col = "I"; Value = 9698.166999999999 },
    @{ Row = 136; Col = "K"; Value = 29094.501 },
    @{ Row = 136; Col = "M"; Value = -26544.501 }
)
foreach ($u in $WVR_updates) {
    $ws.Range($u.Col + $u.Row).Value = $u.Value
}
